$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up cell styles for the new row (A11 should match the rest of the
# row's style, D11 should match the "date text" style used by D7:D10) ---
$ws.Range("B11").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("D10").Copy()
$ws.Range("D11").PasteSpecial(-4122)  # xlPasteFormats

$ws.Application.CutCopyMode = $false | Out-Null

# --- New company data: Unacademy (Sorting Hat Technologies Pvt. Ltd.) ---
$ws.Range("A11").Value = "Sorting Hat Technologies Pvt. Ltd."
$ws.Range("B11").Value = "Unacademy"
$ws.Range("C11").Value = "U72200KA2015PTC082063"
$ws.Range("E11").Value = "Private"
$ws.Range("F11").Value = "Operating"
$ws.Range("G11").Value = "Active"
$ws.Range("J11").Value = "Education Services"
$ws.Range("P11").Value = 8585858585
$ws.Range("S11").Value = "https://unacademy.com/"
$ws.Range("V11").Value = "unacademy.png"
$ws.Range("D11").Value = "6 Aug 2015"

# Row 11 grows to fit the wrapped text, like the other data rows.
$ws.Rows("11:11").RowHeight = 29.4

# Hyperlink for the new website cell. Adding a hyperlink stamps its own
# font style onto the cell, so restore the original "website" cell look
# (border + vertically centred, no special hyperlink font) afterwards,
# matching how S11 looks in the rest of the sheet.
$ws.Hyperlinks.Add($ws.Range("S11"), "https://unacademy.com/")
$ws.Range("T10").Copy()
$ws.Range("S11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false | Out-Null

# Move the active selection to D12, like after typing into D11 and hitting Enter.
$ws.Range("D12").Select() | Out-Null

Write-Host "done"
